# TC06_Canine_Filter_StageOfDisease-4a.xlsx
# Fix the "CasesTab" Cypher query (cell B2 on the "startup" sheet): the
# query was erroneously returning a `Cohort` column (joined via an
# OPTIONAL MATCH on (:cohort) that isn't relevant to this stage-of-disease
# test case). Drop that trailing RETURN item - and the now-dangling comma
# on the previous line - so the query only returns the intended columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$fixedCasesQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)`n" + `
    "MATCH (c)<--(diag:diagnosis)`n" + `
    "OPTIONAL MATCH (samp:sample)-->(c)`n" + `
    "OPTIONAL MATCH (co:cohort)<-[*]-(c)`n" + `
    "WITH DISTINCT c, s, demo, diag, co`n" + `
    "WHERE diag.stage_of_disease IN ['IVa']`n" + `
    "RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,`n" + `
    "        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,`n" + `
    "        coalesce(s.clinical_study_type, '') AS  ``Study Type``,`n" + `
    "        coalesce(demo.breed, '') AS Breed ,`n" + `
    "        coalesce(diag.disease_term, '') AS Diagnosis ,`n" + `
    "        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,`n" + `
    "        coalesce(demo.patient_age_at_enrollment, '') AS Age ,`n" + `
    "        coalesce(demo.sex, '') AS Sex ,`n" + `
    "        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,`n" + `
    "        coalesce(demo.weight, '') AS ``Weight (kg)``,`n" + `
    "        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $fixedCasesQuery

# Leave the selection on the cell that was just edited (matches the
# author's saved cursor position after making the fix).
$ws.Range("B2").Select()
